# feat(phase4): finalize 250 Yemen dataset flow and phase status docs
#
# Mark Phase 4 (WBS rows 63-79) as fully Done: flip the Execution Status
# from "Partial" to "Done", stamp a Completed On date (matching the
# existing Started On date), and fill in the five sign-off columns
# (Schema / Validation / Permissions-Isolation / Workflow / Evidence)
# with check marks, mirroring the pattern already used by the earlier
# completed phases further up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

$checkMark = [char]0x2705

for ($row = 63; $row -le 79; $row++) {
    $ws.Range("H$row").Value = "Done"

    # Completed On (K) should read the same way the existing Started On
    # (J) cell does — copy J's cell (value + style) onto K so it stays a
    # plain text date stamp instead of being reinterpreted as a serial date.
    $ws.Range("J$row").Copy($ws.Range("K$row"))

    $ws.Range("L$row").Value = $checkMark
    $ws.Range("M$row").Value = $checkMark
    $ws.Range("N$row").Value = $checkMark
    $ws.Range("O$row").Value = $checkMark
    $ws.Range("P$row").Value = $checkMark
}

Write-Output "Phase 4 rows 63-79 marked Done with sign-off checks"
